$d = $word.ActiveDocument

# --- Change 1: fix paragraph-mark run properties on the "Vzhladom na to..." paragraph ---
# (remove bold, add w:cs="Times New Roman" to rFonts in the paragraph mark's rPr)
$targetText = "Vzhladom na to"
$p11 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*nepotrebujeme*") {
        $p11 = $p
        break
    }
}
if ($p11 -eq $null) {
    throw "Could not locate target paragraph for formatting fix"
}

$p11Xml = '<w:p w14:paraId="0FBADA58" w14:textId="41AC0E7A" w:rsidR="000A6A51" w:rsidRPr="0007617D" w:rsidRDefault="000A6A51" w:rsidP="00F02DAF"><w:pPr><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cs="Times New Roman"/></w:rPr></w:pPr><w:r><w:tab/><w:t xml:space="preserve">Vzhľadom na to, že nepotrebujeme naše dáta </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>štruktúrovať</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> do nadradených a podradených premenných – objektov ( </w:t></w:r><w:r w:rsidRPr="00F02DAF"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:highlight w:val="lightGray"/></w:rPr><w:t>{“</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00F02DAF"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:highlight w:val="lightGray"/></w:rPr><w:t>sat</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00F02DAF"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:highlight w:val="lightGray"/></w:rPr><w:t>“: {</w:t></w:r><w:r w:rsidR="00374E1D" w:rsidRPr="00F02DAF"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:highlight w:val="lightGray"/></w:rPr><w:t>“GPS“: {“</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00374E1D" w:rsidRPr="00F02DAF"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:highlight w:val="lightGray"/></w:rPr><w:t>lat</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00374E1D" w:rsidRPr="00F02DAF"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:highlight w:val="lightGray"/></w:rPr><w:t>“: 48.582563, “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00374E1D" w:rsidRPr="00F02DAF"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:highlight w:val="lightGray"/></w:rPr><w:t>long</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00374E1D" w:rsidRPr="00F02DAF"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:highlight w:val="lightGray"/></w:rPr><w:t>“: 17.816788}, “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00374E1D" w:rsidRPr="00F02DAF"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:highlight w:val="lightGray"/></w:rPr><w:t>Sens</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00374E1D" w:rsidRPr="00F02DAF"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:highlight w:val="lightGray"/></w:rPr><w:t>“: {</w:t></w:r><w:r w:rsidR="00F02DAF" w:rsidRPr="00F02DAF"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:highlight w:val="lightGray"/></w:rPr><w:t>“T“: 25}}}</w:t></w:r><w:r w:rsidR="00F02DAF"><w:t xml:space="preserve"> ) nám bude stačiť zápis bez objektových/vlnitých zátvoriek. Odstránením úvodzoviek z názvov premenných sa ušetria 2B/premennú. Náš prvotný koncept potom vyzeral takto: </w:t></w:r><w:r w:rsidR="00F02DAF" w:rsidRPr="0007617D"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:highlight w:val="lightGray"/></w:rPr><w:t xml:space="preserve">C: 1, T: 25, H: 40, </w:t></w:r><w:r w:rsidR="008103FD"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:highlight w:val="lightGray"/></w:rPr><w:t xml:space="preserve">P: 1000, </w:t></w:r><w:r w:rsidR="00F02DAF" w:rsidRPr="0007617D"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:highlight w:val="lightGray"/></w:rPr><w:t xml:space="preserve">D: 20, R: 1, A: </w:t></w:r><w:r w:rsidR="0007617D" w:rsidRPr="0007617D"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:highlight w:val="lightGray"/></w:rPr><w:t>48.582563</w:t></w:r><w:r w:rsidR="00F02DAF" w:rsidRPr="0007617D"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:highlight w:val="lightGray"/></w:rPr><w:t>, O: 17.816788</w:t></w:r><w:r w:rsidR="00F02DAF"><w:t>.</w:t></w:r><w:r w:rsidR="0007617D"><w:t xml:space="preserve"> Bolo by možné ešte odstrániť 3 znaky / premennú (dvojbodku, medzeru a čiarku). Dáta by mohli potom vyzerať nasledovne: </w:t></w:r><w:r w:rsidR="0007617D" w:rsidRPr="0007617D"><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t>1C25T40H</w:t></w:r><w:r w:rsidR="008103FD"><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t>1000P</w:t></w:r><w:r w:rsidR="0007617D" w:rsidRPr="0007617D"><w:rPr><w:highlight w:val="lightGray"/></w:rPr><w:t>20D1R</w:t></w:r><w:r w:rsidR="0007617D" w:rsidRPr="0007617D"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:highlight w:val="lightGray"/></w:rPr><w:t>48.582563A17.816788O</w:t></w:r><w:r w:rsidR="0007617D"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="0007617D"><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cs="Times New Roman"/></w:rPr><w:t xml:space="preserve">. </w:t></w:r></w:p>'
$p11.Range.InsertXML($p11Xml)

# --- Change 2: append the new section (page break, heading, and three paragraphs) ---
$endRange = $d.Content
$endRange.Collapse(0)
$newParagraphsXml = '<w:p><w:pPr><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cs="Times New Roman"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cs="Times New Roman"/></w:rPr><w:br w:type="page"/></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Nadpis3"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Vlastný formát – prerábka</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:tab/><w:t>Po konzultáciách sme sa rozhodli ustáliť formát na nasledujúcej štruktúre:</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cs="Times New Roman"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:highlight w:val="lightGray"/></w:rPr><w:t>#123T25.55H48.12P998D10</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cs="Times New Roman"/></w:rPr><w:t xml:space="preserve">pričom: # identifikuje číslo správy, T teplotu, H vlhkosť, P tlak, D vzdialenosť. Do budúcna je možné </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cs="Times New Roman"/></w:rPr><w:t>packet</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cs="Times New Roman"/></w:rPr><w:t xml:space="preserve"> rozšíriť o ďalšie premenné, ako napríklad radiáciu R...</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cs="Times New Roman"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cs="Times New Roman"/></w:rPr><w:tab/></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cs="Times New Roman"/></w:rPr><w:t>Wrapper</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cs="Times New Roman"/></w:rPr><w:t xml:space="preserve"> zostavujúci a rozoberajúci správu z a do jednotlivých premenných by mal byť flexibilný a nemal by závisieť od počtu cifier číselnej hodnoty, alebo poradia jednotlivých premenných v správe.</w:t></w:r></w:p>'
$endRange.InsertXML($newParagraphsXml)

Write-Host "Edit complete. Paragraph count now: $($d.Paragraphs.Count)"
